# Final review ppt.pptx - "Added PPT, Report and Published Paper in DOCUMENTS folder"
#
# The only substantive textual change in this revision is on slide 1
# (the title slide): the subtitle placeholder's text "Batch Nuber: CSE-65"
# is corrected to "Group Number: CSE-65".
#
# The subtitle is made up of several runs (so the "CSE-65" portion keeps
# its own run/formatting). We only touch the two runs that hold the
# misspelled words, rewriting their text in place via Characters() so the
# surrounding run-level formatting (font, language, etc.) is preserved
# exactly as authored.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)                 # "Batch Nuber: CSE-65" subtitle box
$tr = $sh.TextFrame.TextRange

# Run 1: "Batch " (6 chars) -> "Group"
$tr.Characters(1, 6).Text = "Group"

# Run 2 (now starting right after "Group"): "Nuber" (5 chars) -> " Number"
$tr.Characters(6, 5).Text = " Number"
